$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert the new "Roles" slide at position 2 (Title + Content layout).
# ---------------------------------------------------------------------------
$roles = $p.Slides.Add(2, 16)
$roles.Shapes.Item(1).TextFrame.TextRange.Text = "Roles"

$rolesBody = $roles.Shapes.Item(2).TextFrame.TextRange
$rolesBody.Text = "Roberto Ito`rDesigner`rDeveloper`rQA Tester`r`rYayun Yang`rDesigner`rDeveloper`rQA Tester"
$rolesBody.Paragraphs(2, 1).IndentLevel = 2
$rolesBody.Paragraphs(3, 1).IndentLevel = 2
$rolesBody.Paragraphs(4, 1).IndentLevel = 2
$rolesBody.Paragraphs(5, 1).IndentLevel = 2
$rolesBody.Paragraphs(7, 1).IndentLevel = 2
$rolesBody.Paragraphs(8, 1).IndentLevel = 2
$rolesBody.Paragraphs(9, 1).IndentLevel = 2

# ---------------------------------------------------------------------------
# 2. Insert three new slides before the closing "Game link" slide:
#    State Diagram, Class Diagram, Classes (Title and body layout).
# ---------------------------------------------------------------------------
$stateDiagram = $p.Slides.Add(10, 17)
$stateDiagram.Shapes.Item(1).TextFrame.TextRange.Text = "State Diagram"
$stateDiagram.Shapes.Item(1).TextFrame.AutoSize = 2

$classDiagram = $p.Slides.Add(11, 17)
$classDiagram.Shapes.Item(1).TextFrame.TextRange.Text = "Class Diagram"
$classDiagram.Shapes.Item(1).TextFrame.AutoSize = 2

$classes = $p.Slides.Add(12, 17)
$classes.Shapes.Item(1).TextFrame.TextRange.Text = "Classes"
$classes.Shapes.Item(1).TextFrame.AutoSize = 2
